$wb = $excel.ActiveWorkbook

# --- Summary sheet updates ---
$summary = $wb.Worksheets.Item("Summary")

# Update row 5 figures (was 300, now 16.79)
$summary.Range("A5").Value = 16.79
$summary.Range("E5").Value = 16.79
$summary.Range("F5").Value = 16.79

# Row 6 (all zeros) is no longer part of the table - remove it entirely
$summary.Range("A6:F6").Delete()

# --- Repayment schedule sheet updates ---
$repay = $wb.Worksheets.Item("Repayment schedule")

# Row 2: the trailing blank P2 cell is cleared away (O2 stays as a blank cell)
$repay.Range("P2").Clear()

# Row 3
$repay.Range("J3").Value = 8.3699999999999992
$repay.Range("K3").Value = 896.09
$repay.Range("O3").Clear()
$repay.Range("P3").Value = 896.09

# Row 4
$repay.Range("J4").Value = 8.42
$repay.Range("K4").Value = 896.14
$repay.Range("O4").Clear()
$repay.Range("P4").Value = 896.14

# Row 5
$repay.Range("J5").Value = 0
$repay.Range("K5").Value = 887.72
$repay.Range("O5").Clear()
$repay.Range("P5").Value = 887.72

# Rows 6-8: drop the now-unused O column cell, totals (P) stay the same
$repay.Range("O6").Clear()
$repay.Range("O7").Clear()
$repay.Range("O8").Clear()

# --- Active tab / selection bookkeeping ---
# Repayment schedule keeps its own remembered selection (set first so it
# doesn't clobber which sheet ends up active/visible)
$null = $repay.Range("G8").Select()

# Summary becomes the active/visible tab (was NewLoanInput)
$summary.Activate()
$null = $summary.Range("C5").Select()
